$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Thbs1", "Tnfrsf11b", "ECs", 3, 1, 22.700661, 68.10198299999999, 0.08615268874617349, 0.08615268874617349, 1, 0.3333333333333333, 0.06861733333333334, 0.205852, 0.01654048691795588, 0.01654048691795588, 1.557658822724, 14.018929404516, 0.001425007421152808, 0.001425007421152808),
  @("ECs", "Thbs1", "Tnfrsf11b", "FAPs", 3, 1, 22.700661, 68.10198299999999, 0.08615268874617349, 0.08615268874617349, 3, 1, 3.776574666666666, 11.329724, 0.9103586635352137, 0.9103586635352137, 85.73074124918799, 771.5766712426918, 0.07842984658693174, 0.07842984658693174),
  @("ECs", "Thbs1", "Tnfrsf11b", "MuSCs", 3, 1, 22.700661, 68.10198299999999, 0.08615268874617349, 0.08615268874617349, 2, 0.6666666666666666, 0.303255, 0.909765, 0.07310084954683041, 0.07310084954683042, 6.884088951554999, 61.956800563995, 0.006297834738088938, 0.006297834738088938),
  @("FAPs", "Thbs1", "Tnfrsf11b", "ECs", 3, 1, 132.5447616666667, 397.634285, 0.5030288587986086, 0.5030288587986087, 1, 0.3333333333333333, 0.06861733333333334, 0.205852, 0.01654048691795588, 0.01654048691795588, 9.094868092868889, 81.85381283581999, 0.008320342258312664, 0.008320342258312666),
  @("FAPs", "Thbs1", "Tnfrsf11b", "FAPs", 3, 1, 132.5447616666667, 397.634285, 0.5030288587986086, 0.5030288587986087, 3, 1, 3.776574666666666, 11.329724, 0.9103586635352137, 0.9103586635352137, 500.5651891097044, 4505.086701987339, 0.457936679615545, 0.4579366796155451),
  @("FAPs", "Thbs1", "Tnfrsf11b", "MuSCs", 3, 1, 132.5447616666667, 397.634285, 0.5030288587986086, 0.5030288587986087, 2, 0.6666666666666666, 0.303255, 0.909765, 0.07310084954683041, 0.07310084954683042, 40.19486169922499, 361.753755293025, 0.03677183692475088, 0.0367718369247509),
  @("MuSCs", "Thbs1", "Tnfrsf11b", "ECs", 3, 1, 41.94534433333333, 125.836033, 0.159189381961201, 0.159189381961201, 1, 0.3333333333333333, 0.06861733333333334, 0.205852, 0.01654048691795588, 0.01654048691795588, 2.878177673901778, 25.903599065116, 0.002633069889806727, 0.002633069889806727),
  @("MuSCs", "Thbs1", "Tnfrsf11b", "FAPs", 3, 1, 41.94534433333333, 125.836033, 0.159189381961201, 0.159189381961201, 3, 1, 3.776574666666666, 11.329724, 0.9103586635352137, 0.9103586635352137, 158.4097247938769, 1425.687523144892, 0.1449194330111956, 0.1449194330111956),
  @("MuSCs", "Thbs1", "Tnfrsf11b", "MuSCs", 3, 1, 41.94534433333333, 125.836033, 0.159189381961201, 0.159189381961201, 2, 0.6666666666666666, 0.303255, 0.909765, 0.07310084954683041, 0.07310084954683042, 12.720135395805, 114.481218562245, 0.01163687906019867, 0.01163687906019868),
  @("Resolving-Mac", "Thbs1", "Tnfrsf11b", "ECs", 3, 1, 66.302588, 198.907764, 0.2516290704940168, 0.2516290704940168, 1, 0.3333333333333333, 0.06861733333333334, 0.205852, 0.01654048691795588, 0.01654048691795588, 4.549506781658667, 40.945561034928, 0.004162067348683685, 0.004162067348683685),
  @("Resolving-Mac", "Thbs1", "Tnfrsf11b", "FAPs", 3, 1, 66.302588, 198.907764, 0.2516290704940168, 0.2516290704940168, 3, 1, 3.776574666666666, 11.329724, 0.9103586635352137, 0.9103586635352137, 250.3966741752373, 2253.570067577135, 0.2290727043215412, 0.2290727043215412),
  @("Resolving-Mac", "Thbs1", "Tnfrsf11b", "MuSCs", 3, 1, 66.302588, 198.907764, 0.2516290704940168, 0.2516290704940168, 2, 0.6666666666666666, 0.303255, 0.909765, 0.07310084954683041, 0.07310084954683042, 20.10659132394, 180.95932191546, 0.01839429882379191, 0.01839429882379191)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = 2 + $i
  $row = $data[$i]
  for ($j = 0; $j -lt $row.Length; $j++) {
    $ws.Cells.Item($r, $j + 1).Value = $row[$j]
  }
}
